# orders_template.xlsx edit
# - Filters out the two example/sample rows (ORD-001, ORD-002), leaving the
#   header row and a single "placeholder/hint" row describing allowed values.
# - Adjusts a handful of column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (second example order) is removed entirely -----------------
$ws.Rows.Item(3).Delete()

# --- Row 2 becomes a "hint" row instead of a concrete example ----------
# Clear the columns that should end up blank.
$ws.Range("A2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()

# Update the text/number values that changed.
$ws.Range("C2").Value = "냉동 또는 냉장 또는 상온"
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 5000
$ws.Range("H2").Value = 15
$ws.Range("Q2").Value = "Y 또는 N"
$ws.Range("R2").Value = "Y 또는 N"

# B2, K2:P2 and S2 stay exactly as they were.

# --- Column width tweaks -------------------------------------------------
# (Excel's ColumnWidth property is expressed in characters of the Normal
# style font and is offset from the raw OOXML "width" attribute by 5/6 of a
# character; subtract that offset so the saved width matches the target.)
$offset = 5/6

$ws.Columns.Item(1).ColumnWidth = 6 - $offset
$ws.Columns.Item(3).ColumnWidth = 16 - $offset
$ws.Columns.Item(4).ColumnWidth = 9 - $offset
$ws.Columns.Item(5).ColumnWidth = 9 - $offset
$ws.Columns.Item(9).ColumnWidth = 5 - $offset
$ws.Columns.Item(10).ColumnWidth = 6 - $offset
$ws.Columns.Item(17).ColumnWidth = 8 - $offset
$ws.Columns.Item(18).ColumnWidth = 8 - $offset
$ws.Columns.Item(19).ColumnWidth = 6 - $offset
